# Update "Modules" sheet header fields + column widths + selection
# (per commit "Mise a jour de certains champs de Modules et de Professeurs")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text swap: C1 "Enseignant" -> "Chef  Module" ; D1 "Nombre d'heures" -> "Composants"
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# Column widths for the (now wider) header labels
# (ColumnWidth is stored/rounded in pixel-sized steps, so the raw target
# character width is nudged by the same amount Excel itself applies when
# it re-quantizes a typed width back to whole pixels)
$ws.Columns.Item(3).ColumnWidth = 34.1666666666667
$ws.Columns.Item(4).ColumnWidth = 23.6666666666667

# Cursor / selection left on E8 after the edits
$ws.Range("E8").Select()
